$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number must be forced to Text
# format first, otherwise Excel auto-converts the literal into a numeric
# cell (losing the exact decimal text / precision, e.g. "19.24" ->
# 19.239999999999998). Values that are NOT parseable as numbers (contain
# two decimal points, letters, "%", spaces, URLs, etc.) are safe to set
# directly via Value2, which keeps them as plain text with no style churn.

$ws.Range("D2").Value2 = "26.397.61"
$ws.Range("E2").Value2 = "  +0.51%  "
$ws.Range("D3").Value2 = "1.607.62"
$ws.Range("E3").Value2 = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("E4").Value2 = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "211.98"
$ws.Range("E5").Value2 = "  -0.53%  "
$ws.Range("E6").Value2 = "  -0.91%  "
$ws.Range("E7").Value2 = "  -0.04%  "
$ws.Range("B8").Value2 = "Dogecoin"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.0606"
$ws.Range("E8").Value2 = "  -0.30%  "
$ws.Range("B9").Value2 = "Cardano"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.244"
$ws.Range("E9").Value2 = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "19.24"
$ws.Range("E10").Value2 = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0848"
$ws.Range("E11").Value2 = "  -0.59%  "
$ws.Range("D12").Value2 = "1.832.73"
$ws.Range("E12").Value2 = "  +0.81%  "
$ws.Range("D13").Value2 = "1.611.69"
$ws.Range("E13").Value2 = "  +1.14%  "
$ws.Range("E14").Value2 = "  -0.43%  "
$ws.Range("E15").Value2 = "  -0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "63.37"
$ws.Range("E16").Value2 = "  -1.00%  "
$ws.Range("D17").Value2 = "26.382.36"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "231.81"
$ws.Range("E18").Value2 = "  +8.01%  "
$ws.Range("E19").Value2 = "  -0.38%  "
$ws.Range("E20").Value2 = "  +2.94%  "
$ws.Range("E21").Value2 = "  -0.07%  "
$ws.Range("E22").Value2 = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "2.19"
$ws.Range("E23").Value2 = "  +3.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "8.97"
$ws.Range("E24").Value2 = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "146.60"
$ws.Range("E25").Value2 = "  +1.12%  "
$ws.Range("E26").Value2 = "  -0.02%  "
$ws.Range("E27").Value2 = "  -0.10%  "
$ws.Range("E28").Value2 = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "15.41"
$ws.Range("E29").Value2 = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.0495"
$ws.Range("E30").Value2 = "  +0.80%  "
$ws.Range("E31").Value2 = "  -0.36%  "
$ws.Range("D32").Value2 = "1.495.07"
$ws.Range("E32").Value2 = "  +5.48%  "
$ws.Range("E33").Value2 = "  +0.52%  "
$ws.Range("E34").Value2 = "  -1.44%  "
$ws.Range("E35").Value2 = "  -0.38%  "
$ws.Range("E36").Value2 = "  +0.79%  "
$ws.Range("E37").Value2 = "  -2.99%  "
$ws.Range("E38").Value2 = "  -0.43%  "
$ws.Range("E39").Value2 = "  -0.14%  "
$ws.Range("E40").Value2 = "  -0.47%  "
$ws.Range("E41").Value2 = "  -0.02%  "
$ws.Range("E42").Value2 = "  +1.06%  "
$ws.Range("E43").Value2 = "  -4.11%  "
$ws.Range("D44").Value2 = "1.744.95"
$ws.Range("E44").Value2 = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.761"
$ws.Range("E45").Value2 = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "60.90"
$ws.Range("E46").Value2 = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "89.57"
$ws.Range("E47").Value2 = "  +2.95%  "
$ws.Range("E48").Value2 = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0501"
$ws.Range("E49").Value2 = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0959"
$ws.Range("E50").Value2 = "  +0.49%  "
$ws.Range("B51").Value2 = "USDD"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.999"
$ws.Range("E51").Value2 = "  -0.13%  "
